# Refresh the scraped crypto price/volume figures (GitHub Actions daily
# update). Price (column D) and Volume(1h) (column E) cells are stored as
# text in the source sheet; for D-column values that look like plain
# numbers we force the cell to Text format first so Excel's COM layer
# doesn't silently coerce the string to a number (which would drop
# significant trailing zeros, e.g. "251.00" -> 251 or "0.0810" -> 0.081).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.804.73"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "2.524.46"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.06"
$ws.Range("E5").Value = "  +4.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.81"
$ws.Range("E6").Value = "  -1.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  -0.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.93"
$ws.Range("E10").Value = "  -1.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.57"
$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("E13").Value = "  -2.24%  "

$ws.Range("D14").Value = "2.912.44"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.24"
$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.458.15"
$ws.Range("E16").Value = "  -2.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.848"
$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").Value = "42.873.19"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.69"
$ws.Range("E20").Value = "  +4.28%  "

$ws.Range("D21").Value = "0.0₃0964"
$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.83"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.00"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.74"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  +4.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.48"
$ws.Range("E29").Value = "  +7.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.29"
$ws.Range("E30").Value = "  +0.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.99"
$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.82"
$ws.Range("E32").Value = "  +1.42%  "

$ws.Range("E33").Value = "  +3.31%  "

$ws.Range("E34").Value = "  +2.12%  "

$ws.Range("E35").Value = "  -0.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0789"
$ws.Range("E36").Value = "  +0.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.61"
$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("E38").Value = "  -1.55%  "

$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.64"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("E41").Value = "  +14.20%  "

$ws.Range("E42").Value = "  +2.37%  "

$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("E44").Value = "  -1.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.30"
$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("D46").Value = "2.021.43"
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.75"
$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.78"
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").Value = "2.766.44"
$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.54"
$ws.Range("E50").Value = "  +2.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.78"
$ws.Range("E51").Value = "  +1.61%  "
